# Insert a new data row at row 97 (pushing the existing rows 97-153 down
# to 98-154, extending the sheet's used range from A1:T153 to A1:T154),
# then populate the new row 97 with the new "Frambuesa" price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 97..153 down to 98..154 (copies formatting/style from the
# row being pushed down, same as Excel's native "Insert Copied/Entire Row").
$ws.Rows.Item(97).Insert()

# Fill in the newly created row 97 with the new record.
$ws.Range("A97").Value = 9
$ws.Range("B97").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C97").Value = "Metropolitana"
$ws.Range("D97").Value = 45006
$ws.Range("E97").Value = 13
$ws.Range("F97").Value = "Fruta"
$ws.Range("G97").Value = 100101
$ws.Range("H97").Value = "Berries"
$ws.Range("I97").Value = 100101004
$ws.Range("J97").Value = "Frambuesa"
$ws.Range("K97").Value = "Sin especificar"
$ws.Range("L97").Value = "Primera"
$ws.Range("M97").Value = 400
$ws.Range("N97").Value = 6000
$ws.Range("O97").Value = 6500
$ws.Range("P97").Value = 6250
$ws.Range("Q97").Value = "`$/bandeja 2 kilos"
$ws.Range("R97").Value = "Provincia de Linares"
$ws.Range("S97").Value = 3125
$ws.Range("T97").Value = 2
